$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update noOfOrders values (column B)
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("B7").Value = 4

# Update the selected/active cell to C3
$ws.Range("C3").Select()
